$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.036.66'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.186.08'
$ws.Range('E3').Value = '  -2.82%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.33'
$ws.Range('E5').Value = '  -2.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.55'
$ws.Range('E6').Value = '  -6.31%  '
$ws.Range('E7').Value = '  -6.14%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.187.12'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('E10').Value = '  -2.91%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('E12').Value = '  -4.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.737.73'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.139.99'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.47'
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000160'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.150.17'
$ws.Range('E18').Value = '  -2.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '416.42'
$ws.Range('E19').Value = '  -3.04%  '
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.82'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.16'
$ws.Range('E22').Value = '  -3.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.95'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('E25').Value = '  +3.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.493'
$ws.Range('E26').Value = '  -3.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000109'
$ws.Range('E27').Value = '  -3.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.80'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -6.16%  '
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.02'
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.40'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '156.04'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('E37').Value = '  -2.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.755.47'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.06'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('E41').Value = '  -3.08%  '
$ws.Range('E42').Value = '  -6.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.94'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.72'
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0629'
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.10'
$ws.Range('E46').Value = '  -4.51%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '298.70'
$ws.Range('E47').Value = '  -4.99%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0264'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.08'
$ws.Range('E49').Value = '  -9.57%  '
$ws.Range('E50').Value = '  -5.85%  '
$ws.Range('E51').Value = '  -0.05%  '
